$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Global Annual Qty multiplier: 1000 -> 1
$ws.Range("L1").Value = 1

# 2. New row 13: USB C Receptacle (moved up from old row 28)
$ws.Range("B13").Value = "USB C Recepticle"
$ws.Range("C13").Value = 1
$ws.Range("D13").Value = "USB4105-GF-A"
$ws.Range("E13").Value = "Connector"
$ws.Range("M13").Value = 0.78

# 3. New row 14: 5.1k Resistor (moved up from old row 29)
$ws.Range("B14").Value = "5.1k Resistor"
$ws.Range("C14").Value = 2
$ws.Range("D14").Value = "RC0603FR-075K1L"
$ws.Range("E14").Value = "Resistor"
$ws.Range("F14").Value = "5.1k"
$ws.Range("K14").Value = "0603"
$ws.Range("M14").Value = 0.1

# 4. New row 15: Temperature Sensor (moved up from old row 31, part number updated)
$ws.Range("B15").Value = "Temperature Sensor"
$ws.Range("C15").Value = 1
$ws.Range("D15").Value = "MCP9700AT-E/TT"
$ws.Range("E15").Value = "Sensor"
$ws.Range("K15").Value = "SOT"
$ws.Range("M15").Value = 0.48

# 5. New row 16: 10k Resistor (brand new pushbutton pull-up resistor)
$ws.Range("B16").Value = "10k Resistor"
$ws.Range("C16").Value = 1
$ws.Range("D16").Value = "RMCF0603JT10K0"
$ws.Range("E16").Value = "Resistor"
$ws.Range("F16").Value = "10k"
$ws.Range("K16").Value = "0603"
$ws.Range("M16").Value = 0.1

# 6. Clear old row 28 (USB C Receptacle data moved to row 13)
$ws.Range("B28").Value = ""
$ws.Range("C28").Value = ""
$ws.Range("D28").Value = ""
$ws.Range("E28").Value = ""
$ws.Range("M28").Value = ""

# 7. Clear old row 29 (5.1k Resistor data moved to row 14)
$ws.Range("B29").Value = ""
$ws.Range("C29").Value = ""
$ws.Range("D29").Value = ""
$ws.Range("E29").Value = ""
$ws.Range("F29").Value = ""
$ws.Range("I29").Value = ""
$ws.Range("K29").Value = ""
$ws.Range("M29").Value = ""

# 8. Clear old row 31 (Temperature Sensor data moved to row 15)
$ws.Range("B31").Value = ""
$ws.Range("C31").Value = ""
$ws.Range("D31").Value = ""
$ws.Range("E31").Value = ""
$ws.Range("K31").Value = ""
$ws.Range("M31").Value = ""

# 9. Update selection to match final cursor position
$ws.Range("E18").Select()
